$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch rows 11-23 so they materialize as empty row stubs in the sheet,
# mirroring the rows that Excel left behind while the Move page handled
# a null-location move (no data in these rows).
for ($r = 11; $r -le 23; $r++) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# New report row (row 24) describing the move-with-null-location event.
$newRow = 24
$ws.Cells.Item($newRow, 1).Value = 23
$ws.Cells.Item($newRow, 2).Value = "2025-04-26 20:43:43"
$ws.Cells.Item($newRow, 3).Value = "John Smith moved battery 7 from No Location to floor space 1.`nNow John Smith is Confident.`n"

# Match the wrap-text styling used by the other report rows in column C.
$ws.Cells.Item($newRow, 3).WrapText = $true

# Undo the explicit row-height override that wrapping the long text would
# otherwise leave behind, so the row keeps using the sheet's default height.
$ws.Rows.Item($newRow).AutoFit()
